# Especificaciones.xlsx — add "Estabilidad LDO" worksheet with worst-case
# voltage-loop stability notes, and move the active view/selection onto it.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet, inserted right after "Sheet1".
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1, [System.Type]::Missing, [System.Type]::Missing)
$ws2.Name = "Estabilidad LDO"

# Content (order matches the shared-string build order of the target file:
# A6 is typed before A5).
$ws2.Range("A1").Value = "Peor caso lazo tensión"
$ws2.Range("A3").Value = "RL = 250"
$ws2.Range("A4").Value = "CL = 15uF"
$ws2.Range("A6").Value = "P2 = 22.8kHz"
$ws2.Range("A5").Value = "Queda:"

# Leave the old sheet's selection parked further down the sheet ...
$ws1.Range("A48").Select() | Out-Null

# ... and finish with the new sheet active/selected, cursor on the last line.
$ws2.Activate() | Out-Null
$ws2.Range("A6").Select() | Out-Null
